# Auto-generated edit script: updates computed profit/price columns (H-N)
# on each job sheet's leve table, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4538.76
$ws.Range("I32").Value = 4620
$ws.Range("K32").Value = 4620
$ws.Range("M32").Value = -4294
$ws.Range("H33").Value = 41667140
$ws.Range("I33").Value = 303.3
$ws.Range("K33").Value = 303.3
$ws.Range("M33").Value = -74.30000000000001
$ws.Range("H38").Value = 2105.4443
$ws.Range("I38").Value = 868.9091
$ws.Range("J38").Value = 4048.5715
$ws.Range("K38").Value = 2606.7273
$ws.Range("L38").Value = 12145.7145
$ws.Range("M38").Value = -2234.7273
$ws.Range("N38").Value = -12889.7145
$ws.Range("H41").Value = 256
$ws.Range("I41").Value = 249.5
$ws.Range("K41").Value = 249.5
$ws.Range("M41").Value = 190.5
$ws.Range("H42").Value = 2312.8572
$ws.Range("I42").Value = 25.25
$ws.Range("J42").Value = 5363
$ws.Range("K42").Value = 75.75
$ws.Range("L42").Value = 16089
$ws.Range("M42").Value = 154.25
$ws.Range("N42").Value = -16549
$ws.Range("H45").Value = 5850
$ws.Range("I45").Value = 9000
$ws.Range("K45").Value = 27000
$ws.Range("M45").Value = -26808
$ws.Range("H52").Value = 2017.8572
$ws.Range("I52").Value = 781.75
$ws.Range("J52").Value = 3666
$ws.Range("K52").Value = 2345.25
$ws.Range("L52").Value = 10998
$ws.Range("M52").Value = -2185.25
$ws.Range("N52").Value = -11318
$ws.Range("H132").Value = 1741.1666
$ws.Range("I132").Value = 1623.92
$ws.Range("K132").Value = 4871.76
$ws.Range("M132").Value = -2341.76

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3967.0952
$ws.Range("I61").Value = 2612.353
$ws.Range("K61").Value = 2612.353
$ws.Range("M61").Value = -2400.353
$ws.Range("H102").Value = 31657.227
$ws.Range("I102").Value = 11065.467
$ws.Range("J102").Value = 75782.42999999999
$ws.Range("K102").Value = 11065.467
$ws.Range("L102").Value = 75782.42999999999
$ws.Range("M102").Value = -9443.467000000001
$ws.Range("N102").Value = -79026.42999999999
$ws.Range("H131").Value = 73495
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 73495
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 73495
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -83575
$ws.Range("H136").Value = 3967.0952
$ws.Range("I136").Value = 2612.353
$ws.Range("K136").Value = 7837.059
$ws.Range("M136").Value = -5287.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1026.1852
$ws.Range("I94").Value = 541.3913
$ws.Range("J94").Value = 3813.75
$ws.Range("K94").Value = 541.3913
$ws.Range("L94").Value = 3813.75
$ws.Range("M94").Value = -90.3913
$ws.Range("N94").Value = -4715.75
$ws.Range("H134").Value = 961549.5600000001
$ws.Range("I134").Value = 1318755.6
$ws.Range("K134").Value = 3956266.8
$ws.Range("M134").Value = -3953731.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3644.6428
$ws.Range("I31").Value = 3252.7
$ws.Range("J31").Value = 4624.5
$ws.Range("K31").Value = 3252.7
$ws.Range("L31").Value = 4624.5
$ws.Range("M31").Value = -2957.7
$ws.Range("N31").Value = -5214.5
$ws.Range("H34").Value = 3644.6428
$ws.Range("I34").Value = 3252.7
$ws.Range("J34").Value = 4624.5
$ws.Range("K34").Value = 3252.7
$ws.Range("L34").Value = 4624.5
$ws.Range("M34").Value = -3050.7
$ws.Range("N34").Value = -5028.5
$ws.Range("H132").Value = 3589.3572
$ws.Range("I132").Value = 3380.08
$ws.Range("J132").Value = 5333.3335
$ws.Range("K132").Value = 10140.24
$ws.Range("L132").Value = 16000.0005
$ws.Range("M132").Value = -7610.24
$ws.Range("N132").Value = -21060.0005
$ws.Range("H134").Value = 2529.7585
$ws.Range("I134").Value = 2311.95
$ws.Range("J134").Value = 3013.7778
$ws.Range("K134").Value = 6935.849999999999
$ws.Range("L134").Value = 9041.3334
$ws.Range("M134").Value = -4400.849999999999
$ws.Range("N134").Value = -14111.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 70320.664
$ws.Range("I114").Value = 900
$ws.Range("J114").Value = 105031
$ws.Range("K114").Value = 2700
$ws.Range("L114").Value = 315093
$ws.Range("M114").Value = 554
$ws.Range("N114").Value = -321601
$ws.Range("H117").Value = 1000031
$ws.Range("J117").Value = 1000031
$ws.Range("L117").Value = 3000093
$ws.Range("N117").Value = -3006977
$ws.Range("H120").Value = 10166.5
$ws.Range("I120").Value = 7249.75
$ws.Range("J120").Value = 16000
$ws.Range("K120").Value = 21749.25
$ws.Range("L120").Value = 48000
$ws.Range("M120").Value = -16911.25
$ws.Range("N120").Value = -57676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 15069.6
$ws.Range("I55").Value = 350
$ws.Range("J55").Value = 18749.5
$ws.Range("K55").Value = 350
$ws.Range("L55").Value = 18749.5
$ws.Range("M55").Value = -23
$ws.Range("N55").Value = -19403.5
$ws.Range("H98").Value = 64300
$ws.Range("J98").Value = 64300
$ws.Range("L98").Value = 64300
$ws.Range("N98").Value = -70290
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H104").Value = 66050
$ws.Range("J104").Value = 66050
$ws.Range("L104").Value = 66050
$ws.Range("N104").Value = -73038
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H110").Value = 69575.375
$ws.Range("J110").Value = 69575.375
$ws.Range("L110").Value = 69575.375
$ws.Range("N110").Value = -77755.375
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H118").Value = 100000
$ws.Range("J118").Value = 100000
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6907.625
$ws.Range("I40").Value = 5877
$ws.Range("K40").Value = 5877
$ws.Range("M40").Value = -5741
$ws.Range("H61").Value = 3060.8
$ws.Range("I61").Value = 3422.125
$ws.Range("J61").Value = 1615.5
$ws.Range("K61").Value = 3422.125
$ws.Range("L61").Value = 1615.5
$ws.Range("M61").Value = -3220.125
$ws.Range("N61").Value = -2019.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 22222
$ws.Range("J106").Value = 22222
$ws.Range("L106").Value = 22222
$ws.Range("N106").Value = -24746
$ws.Range("H113").Value = 3060.8
$ws.Range("I113").Value = 3422.125
$ws.Range("J113").Value = 1615.5
$ws.Range("K113").Value = 3422.125
$ws.Range("L113").Value = 1615.5
$ws.Range("M113").Value = -1252.125
$ws.Range("N113").Value = -5955.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 98626.336
$ws.Range("J116").Value = 98626.336
$ws.Range("L116").Value = 98626.336
$ws.Range("N116").Value = -107804.336
$ws.Range("H117").Value = 61499
$ws.Range("J117").Value = 61499
$ws.Range("L117").Value = 61499
$ws.Range("N117").Value = -70677
$ws.Range("H119").Value = 57500
$ws.Range("J119").Value = 57500
$ws.Range("L119").Value = 57500
$ws.Range("N119").Value = -67176
$ws.Range("H128").Value = 100858
$ws.Range("J128").Value = 100858
$ws.Range("L128").Value = 100858
$ws.Range("N128").Value = -110818
$ws.Range("H132").Value = 1021.0417
$ws.Range("I132").Value = 1020.4545
$ws.Range("J132").Value = 1027.5
$ws.Range("K132").Value = 3061.3635
$ws.Range("L132").Value = 3082.5
$ws.Range("M132").Value = -531.3635000000004
$ws.Range("N132").Value = -8142.5
$ws.Range("H136").Value = 4641.3335
$ws.Range("I136").Value = 4646.9443
$ws.Range("J136").Value = 4624.5
$ws.Range("K136").Value = 13940.8329
$ws.Range("L136").Value = 13873.5
$ws.Range("M136").Value = -11390.8329
$ws.Range("N136").Value = -18973.5
